# Update the "dSF" (column F) values to match repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value  = -1
$ws.Range("F7").Value  = 3
$ws.Range("F11").Value = -5
$ws.Range("F13").Value = -6
$ws.Range("F14").Value = -2
$ws.Range("F17").Value = -1
$ws.Range("F20").Value = -4
$ws.Range("F22").Value = -4
$ws.Range("F25").Value = -1
$ws.Range("F26").Value = -3
$ws.Range("F27").Value = -2
$ws.Range("F31").Value = 1
